$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels to append " %"
$ws.Range("E1").Value = "Spain Monthly Inflation Rate %"
$ws.Range("J1").Value = "US Monthly Inflation Rate %"

# Flip the sign of every numeric value in column L (Real Exchange Rate Growth),
# for data rows 3 through 373.
for ($r = 3; $r -le 373; $r++) {
    $cell = $ws.Cells.Item($r, 12)
    $val = $cell.Value()
    if ($val -ne $null) {
        $cell.Value = -1 * $val
    }
}
